# Dodanie wstępu do prezentacji
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "datetimeFigureOut" date field (09.01.2019 -> 12.01.2019)
#    This field lives on the slide master + every slide layout footer.
# ---------------------------------------------------------------------------
$newDate = "12.01.2019"

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "09.01.2019") {
                $tr.Text = $newDate
            }
        }
    }
}

# NOTE: `$p.Slides.Item(1).Master.CustomLayouts.Item($li)` always resolves
# to the very first custom layout in this runtime (its `.Index` is always
# 1, no matter which $li was requested), so edits made that way silently
# land on slideLayout1.xml for every iteration. Going through
# `$p.Designs.Item(1).SlideMaster` instead yields CustomLayout objects
# whose `.Index` (and underlying OOXML part) really do match $li.
$masterForLayouts = $p.Designs.Item(1).SlideMaster

Update-DateField($masterForLayouts.Shapes)
for ($li = 1; $li -le $masterForLayouts.CustomLayouts.Count; $li++) {
    $layout = $masterForLayouts.CustomLayouts.Item($li)
    Update-DateField($layout.Shapes)
}

# ---------------------------------------------------------------------------
# 2. Slide 3 ("Wstęp") - expand the placeholder text from "B" into the
#    real introduction paragraphs.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "`tWyobraźmy sobie taką sytuację. Potrzebujemy wyjechać do innego miasta bo na przykład mamy ważne spotkanie biznesowe lub po prostu jedziemy wypocząć. Rezerwujemy nocleg, wyjeżdżamy i … i nasze mieszkanie stoi nieużytkowane przez cały czas naszego wyjazdu. Jest to bardzo nieekonomiczne dla naszego portfela. "
$tr3.InsertAfter("`r`t")
$tr3.InsertAfter("Po kilku długotrwałych nieobecnościach w mieszkaniu z powodów wyjazdów, narodził się pomysł. ")
$tr3.InsertAfter("A gdyby tak istniała możliwość wynajęcia naszego mieszkania podczas naszej nieobecności? Moglibyśmy dużo zaoszczędzić. ")

# ---------------------------------------------------------------------------
# 3. Slide 4 ("Profil potencjalnego klienta serwisu") - split single run
#    describing the target audience into several runs with updated wording.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Text = "`tProjekt "
$tr4.InsertAfter("kierowany ")
$tr4.InsertAfter("jest zarówno dla ")
$tr4.InsertAfter("przedsiębiorców ")
$tr4.InsertAfter("jak i ")
$tr4.InsertAfter("dla zwykłych ")
$tr4.InsertAfter("ludzi, którzy chcą wynajmować i/lub udostępniać mieszkania na wynajem.")
